# Update the assignment guidelines workbook:
#  - Rebrand "TCSM" course references to "ISEM"
#  - Rename "5. Discussion of results" -> "5. Discussion of findings"
#  - Fix a couple of typos ("innacurate" -> "inaccurate", "Importasnt" -> "Important")
#  - Reword the formatting/submission gross-difference bullet
#  - Leave the active selection on D6, as it was when last saved

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = "ISEM Assignment 1: Path Analysis"

$ws.Cells.Item(11, 1).Value = "5. Discussion of findings"

$ws.Cells.Item(6, 1).Value = "Missing or malformed research question; Insufficient model (e.g., fewer than three variables); Missing or inaccurate path diagram; Theoretical model is not described well enough to interpret the results; Theoretical model does not match the RQ."

$ws.Cells.Item(10, 1).Value = "Results not interpreted or interpreted incorrectly; Applicable measures of explained variance are not included; Important results are omitted; Assumptions are not checked or checked/evaluated incorrectly; Results are included by copying/embedding R output."

$ws.Cells.Item(14, 1).Value = "Gross differences from the required format (e.g., wrong file type, missing sections, submitting redundant files, missing author information)."

$ws.Range("D6").Select()
